$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 39: re-label "RaceData (not a RaceData object)" as "PlayerRace" ---
# D39 already shows this text (shared string reused), so renaming it in place
# (first, before any brand-new strings are introduced) updates the shared
# string for every cell that points at it and keeps its original index.
$ws.Range("D39").Value2 = "PlayerRace"

# K39 used to be struck-through ("RaceData" shown crossed out, s="1").
# It is re-enabled (strikethrough removed) and now shows "PlayerRace" too,
# picking up a new cell style (applyFont but otherwise default formatting).
$ws.Range("K39").Value2 = "PlayerRace"
$ws.Range("K39").Font.Strikethrough = $false
$ws.Range("K39").Font.Name = "Calibri"

# L39 is a newly added cell mirroring the same value, unstyled.
$ws.Range("L39").Value2 = "PlayerRace"

# --- Row 32: BattlePlans -> BattlePlans (in PlayerData) ---
# K32 picks up a new shared string; L32 (duplicate) is removed entirely.
$ws.Range("K32").Value2 = "BattlePlans (in PlayerData)"
$ws.Range("L32").ClearContents()

# --- Row 34: PlayerRelations -> PlayerRelations (in PlayerData) ---
# K34 picks up a new shared string; L34 (duplicate) is removed entirely.
$ws.Range("K34").Value2 = "PlayerRelations (in PlayerData)"
$ws.Range("L34").ClearContents()

# --- Row 47: component editor "Hull map" re-enabled -> K47 now also shows it ---
# K47 is a new cell (struck-through, like the other disabled-field markers
# in this column) showing "PlayerData", matching E47/L47.
$ws.Range("K47").Value2 = "PlayerData"
$ws.Range("K47").Font.Strikethrough = $true

# --- Sheet view: selection moved, scrolled view reset to top ---
$ws.Range("I36").Select() | Out-Null
